$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1157.762
$ws.Range("I32").Value = 98
$ws.Range("J32").Value = 1210.75
$ws.Range("K32").Value = 98
$ws.Range("L32").Value = 1210.75
$ws.Range("M32").Value = 228
$ws.Range("N32").Value = -1862.75

$ws.Range("H64").Value = 6479.8
$ws.Range("I64").Value = 5414.143
$ws.Range("K64").Value = 5414.143
$ws.Range("M64").Value = -5166.143

$ws.Range("H67").Value = 6479.8
$ws.Range("I67").Value = 5414.143
$ws.Range("K67").Value = 5414.143
$ws.Range("M67").Value = -4556.143

$ws.Range("H76").Value = 37516.414
$ws.Range("I76").Value = 39962.074
$ws.Range("K76").Value = 39962.074
$ws.Range("M76").Value = -39647.074

$ws.Range("H79").Value = 37516.414
$ws.Range("I79").Value = 39962.074
$ws.Range("K79").Value = 39962.074
$ws.Range("M79").Value = -38870.074

$ws.Range("H96").Value = 515.6111
$ws.Range("I96").Value = 527.36365
$ws.Range("J96").Value = 497.14285
$ws.Range("K96").Value = 1582.09095
$ws.Range("L96").Value = 1491.42855
$ws.Range("M96").Value = -209.09095
$ws.Range("N96").Value = -4237.428550000001

$ws.Range("H107").Value = 30304704
$ws.Range("I107").Value = 1840.4
$ws.Range("K107").Value = 1840.4
$ws.Range("M107").Value = 79.59999999999991

$ws.Range("H113").Value = 3670.5715
$ws.Range("I113").Value = 3708.3333
$ws.Range("J113").Value = 3444
$ws.Range("K113").Value = 3708.3333
$ws.Range("L113").Value = 3444
$ws.Range("M113").Value = -454.3332999999998
$ws.Range("N113").Value = -9952

$ws.Range("H132").Value = 6945912.5
$ws.Range("J132").Value = 2664.2
$ws.Range("L132").Value = 7992.599999999999
$ws.Range("N132").Value = -13052.6

$ws.Range("H138").Value = 7297.778
$ws.Range("I138").Value = 11168.454
$ws.Range("J138").Value = 6045.5
$ws.Range("K138").Value = 33505.362
$ws.Range("L138").Value = 18136.5
$ws.Range("M138").Value = -28365.362
$ws.Range("N138").Value = -28416.5

$ws.Range("H141").Value = 3667.182
$ws.Range("I141").Value = 2533.9
$ws.Range("K141").Value = 7601.700000000001
$ws.Range("M141").Value = -2421.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 542.3333
$ws.Range("I5").Value = 582.1111
$ws.Range("J5").Value = 423
$ws.Range("K5").Value = 582.1111
$ws.Range("L5").Value = 423
$ws.Range("M5").Value = -470.1111
$ws.Range("N5").Value = -647

$ws.Range("H32").Value = 2088.5217
$ws.Range("I32").Value = 881.7531
$ws.Range("K32").Value = 881.7531
$ws.Range("M32").Value = -594.7531

$ws.Range("H63").Value = 19008.041
$ws.Range("I63").Value = 3481.125
$ws.Range("J63").Value = 26771.5
$ws.Range("K63").Value = 3481.125
$ws.Range("L63").Value = 26771.5
$ws.Range("M63").Value = -2795.125
$ws.Range("N63").Value = -28143.5

$ws.Range("H66").Value = 19008.041
$ws.Range("I66").Value = 3481.125
$ws.Range("J66").Value = 26771.5
$ws.Range("K66").Value = 17405.625
$ws.Range("L66").Value = 133857.5
$ws.Range("M66").Value = -13973.625
$ws.Range("N66").Value = -140721.5

$ws.Range("H74").Value = 623478.5
$ws.Range("I74").Value = 4756.926
$ws.Range("J74").Value = 2479643.2
$ws.Range("K74").Value = 4756.926
$ws.Range("L74").Value = 2479643.2
$ws.Range("M74").Value = -3882.926
$ws.Range("N74").Value = -2481391.2

$ws.Range("H77").Value = 623478.5
$ws.Range("I77").Value = 4756.926
$ws.Range("J77").Value = 2479643.2
$ws.Range("K77").Value = 23784.63
$ws.Range("L77").Value = 12398216
$ws.Range("M77").Value = -19416.63
$ws.Range("N77").Value = -12406952

$ws.Range("H94").Value = 37000
$ws.Range("J94").Value = 37000
$ws.Range("L94").Value = 37000
$ws.Range("N94").Value = -38802

$ws.Range("H102").Value = 2135.0952
$ws.Range("I102").Value = 2356.4
$ws.Range("K102").Value = 2356.4
$ws.Range("M102").Value = -734.4000000000001

$ws.Range("H122").Value = 6645.8887
$ws.Range("I122").Value = 6971.8335
$ws.Range("K122").Value = 20915.5005
$ws.Range("M122").Value = -18465.5005

$ws.Range("H134").Value = 40467.11
$ws.Range("J134").Value = 40467.11
$ws.Range("L134").Value = 40467.11
$ws.Range("N134").Value = -50607.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 542.3333
$ws.Range("I4").Value = 582.1111
$ws.Range("J4").Value = 423
$ws.Range("K4").Value = 582.1111
$ws.Range("L4").Value = 423
$ws.Range("M4").Value = -467.1111
$ws.Range("N4").Value = -653

$ws.Range("H105").Value = 14486.5
$ws.Range("I105").Value = 14735.111
$ws.Range("J105").Value = 14283.091
$ws.Range("K105").Value = 14735.111
$ws.Range("L105").Value = 14283.091
$ws.Range("M105").Value = -12988.111
$ws.Range("N105").Value = -17777.091

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4143.514
$ws.Range("I31").Value = 3874.5715
$ws.Range("J31").Value = 4322.8096
$ws.Range("K31").Value = 3874.5715
$ws.Range("L31").Value = 4322.8096
$ws.Range("M31").Value = -3579.5715
$ws.Range("N31").Value = -4912.8096

$ws.Range("H34").Value = 4143.514
$ws.Range("I34").Value = 3874.5715
$ws.Range("J34").Value = 4322.8096
$ws.Range("K34").Value = 3874.5715
$ws.Range("L34").Value = 4322.8096
$ws.Range("M34").Value = -3672.5715
$ws.Range("N34").Value = -4726.8096

$ws.Range("H64").Value = 48999.5
$ws.Range("J64").Value = 48999.5
$ws.Range("L64").Value = 48999.5
$ws.Range("N64").Value = -49495.5

$ws.Range("H67").Value = 48999.5
$ws.Range("J67").Value = 48999.5
$ws.Range("L67").Value = 48999.5
$ws.Range("N67").Value = -50715.5

$ws.Range("H86").Value = 7165.3438
$ws.Range("I86").Value = 2951.0454
$ws.Range("K86").Value = 2951.0454
$ws.Range("M86").Value = -1828.0454

$ws.Range("H89").Value = 7165.3438
$ws.Range("I89").Value = 2951.0454
$ws.Range("K89").Value = 14755.227
$ws.Range("M89").Value = -9139.226999999999

$ws.Range("H109").Value = 48199.6
$ws.Range("J109").Value = 48199.6
$ws.Range("L109").Value = 48199.6
$ws.Range("N109").Value = -50279.6

$ws.Range("H132").Value = 6668783.5
$ws.Range("I132").Value = 1923.8485
$ws.Range("J132").Value = 19610336
$ws.Range("K132").Value = 5771.5455
$ws.Range("L132").Value = 58831008
$ws.Range("M132").Value = -3241.5455
$ws.Range("N132").Value = -58836068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1109
$ws.Range("I5").Value = 532.0714
$ws.Range("J5").Value = 1782.0834
$ws.Range("K5").Value = 1596.2142
$ws.Range("L5").Value = 5346.2502
$ws.Range("M5").Value = -1484.2142
$ws.Range("N5").Value = -5570.2502

$ws.Range("H33").Value = 9136902
$ws.Range("I33").Value = 110.44444
$ws.Range("J33").Value = 15462373
$ws.Range("K33").Value = 662.66664
$ws.Range("L33").Value = 92774238
$ws.Range("M33").Value = -379.66664
$ws.Range("N33").Value = -92774804

$ws.Range("H75").Value = 1756.2941
$ws.Range("I75").Value = 805.4286
$ws.Range("J75").Value = 2421.9
$ws.Range("K75").Value = 2416.2858
$ws.Range("L75").Value = 7265.700000000001
$ws.Range("M75").Value = -1418.2858
$ws.Range("N75").Value = -9261.700000000001

$ws.Range("H78").Value = 1756.2941
$ws.Range("I78").Value = 805.4286
$ws.Range("J78").Value = 2421.9
$ws.Range("K78").Value = 7248.8574
$ws.Range("L78").Value = 21797.1
$ws.Range("M78").Value = -2256.8574
$ws.Range("N78").Value = -31781.1

$ws.Range("H107").Value = 666.8570999999999
$ws.Range("J107").Value = 640.6667
$ws.Range("L107").Value = 1922.0001
$ws.Range("N107").Value = -5762.0001

$ws.Range("H114").Value = 1362.4166
$ws.Range("I114").Value = 1083.8334
$ws.Range("K114").Value = 3251.5002
$ws.Range("M114").Value = 2.49980000000005

$ws.Range("H119").Value = 11947.083
$ws.Range("I119").Value = 5922.5
$ws.Range("K119").Value = 17767.5
$ws.Range("M119").Value = -12929.5

$ws.Range("H132").Value = 47225.547
$ws.Range("J132").Value = 78832.16
$ws.Range("L132").Value = 709489.4400000001
$ws.Range("N132").Value = -714549.4400000001

$ws.Range("H135").Value = 1109
$ws.Range("I135").Value = 532.0714
$ws.Range("J135").Value = 1782.0834
$ws.Range("K135").Value = 4788.6426
$ws.Range("L135").Value = 16038.7506
$ws.Range("M135").Value = -2253.6426
$ws.Range("N135").Value = -21108.7506

$ws.Range("H137").Value = 1917.7646
$ws.Range("J137").Value = 2518.4443
$ws.Range("L137").Value = 7555.3329
$ws.Range("N137").Value = -17755.3329

$ws.Range("H140").Value = 1551.0204
$ws.Range("I140").Value = 1106.4681
$ws.Range("K140").Value = 3319.4043
$ws.Range("M140").Value = 1860.5957

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 47338.8
$ws.Range("I62").Value = 45564.668
$ws.Range("K62").Value = 45564.668
$ws.Range("M62").Value = -44878.668

$ws.Range("H65").Value = 47338.8
$ws.Range("I65").Value = 45564.668
$ws.Range("K65").Value = 136694.004
$ws.Range("M65").Value = -133262.004

$ws.Range("H70").Value = 3134.121
$ws.Range("I70").Value = 2876.1482
$ws.Range("K70").Value = 2876.1482
$ws.Range("M70").Value = -2606.1482

$ws.Range("H73").Value = 3134.121
$ws.Range("I73").Value = 2876.1482
$ws.Range("K73").Value = 2876.1482
$ws.Range("M73").Value = -1940.1482

$ws.Range("H94").Value = 33999
$ws.Range("J94").Value = 33999
$ws.Range("L94").Value = 33999
$ws.Range("N94").Value = -35351

$ws.Range("H102").Value = 31252348
$ws.Range("I102").Value = 35715970
$ws.Range("K102").Value = 35715970
$ws.Range("M102").Value = -35714348

$ws.Range("H110").Value = 112211
$ws.Range("J110").Value = 112211
$ws.Range("L110").Value = 112211
$ws.Range("N110").Value = -120391

$ws.Range("H141").Value = 81104.25
$ws.Range("J141").Value = 78472.336
$ws.Range("L141").Value = 78472.336
$ws.Range("N141").Value = -88832.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4472.9355
$ws.Range("I40").Value = 4228.4165
$ws.Range("K40").Value = 4228.4165
$ws.Range("M40").Value = -4092.4165

$ws.Range("H55").Value = 1243.7037
$ws.Range("I55").Value = 1421.125
$ws.Range("J55").Value = 985.63635
$ws.Range("K55").Value = 1421.125
$ws.Range("L55").Value = 985.63635
$ws.Range("M55").Value = -1248.125
$ws.Range("N55").Value = -1331.63635

$ws.Range("H56").Value = 18000
$ws.Range("I56").Value = 18000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 18000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -17309
$ws.Range("N56").ClearContents()

$ws.Range("H61").Value = 2537.6047
$ws.Range("I61").Value = 1607.8055
$ws.Range("K61").Value = 1607.8055
$ws.Range("M61").Value = -1405.8055

$ws.Range("H109").Value = 87000
$ws.Range("J109").Value = 87000
$ws.Range("L109").Value = 87000
$ws.Range("N109").Value = -89774

$ws.Range("H112").Value = 138752.22
$ws.Range("J112").Value = 138752.22
$ws.Range("L112").Value = 138752.22
$ws.Range("N112").Value = -141706.22

$ws.Range("H113").Value = 2537.6047
$ws.Range("I113").Value = 1607.8055
$ws.Range("K113").Value = 1607.8055
$ws.Range("M113").Value = 562.1945000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 52335.75
$ws.Range("J95").Value = 52335.75
$ws.Range("L95").Value = 52335.75
$ws.Range("N95").Value = -57827.75

$ws.Range("H122").Value = 3918.8667
$ws.Range("I122").Value = 3435
$ws.Range("J122").Value = 5249.5
$ws.Range("K122").Value = 10305
$ws.Range("L122").Value = 15748.5
$ws.Range("M122").Value = -7855
$ws.Range("N122").Value = -20648.5

$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840

$ws.Range("H126").Value = 2842.5334
$ws.Range("I126").Value = 2356.7693
$ws.Range("K126").Value = 7070.3079
$ws.Range("M126").Value = -4600.3079

$ws.Range("H132").Value = 1690.3684
$ws.Range("I132").Value = 1019.75
$ws.Range("J132").Value = 3568.1
$ws.Range("K132").Value = 3059.25
$ws.Range("L132").Value = 10704.3
$ws.Range("M132").Value = -529.25
$ws.Range("N132").Value = -15764.3
